$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (dates 2021-04-15, 2021-04-16, 2021-04-17 -> serials 44301-44303)
$newRows = @(
    @{ Row = 227; A = 44301; B = 0; C = 2; D = 43.62050163576881 },
    @{ Row = 228; A = 44302; B = 1; C = 3; D = 65.43075245365321 },
    @{ Row = 229; A = 44303; B = 0; C = 3; D = 65.43075245365321 }
)

foreach ($r in $newRows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D

    # Column A keeps the same date-style formatting as the preceding rows
    $ws.Cells.Item(226, 1).Copy()
    $ws.Cells.Item($r.Row, 1).PasteSpecial(-4122)
}
